$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# ---- Title placeholder ("Title 1", shape id 2) ----
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Data exploration`t"
$titleTr.Font.Language = 1033

# ---- Content placeholder ("Content Placeholder 2", shape id 3) ----
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

# Build the five paragraphs. A one-character placeholder ("X") is kept at the
# very end (in the position of the slide's original, already-empty paragraph)
# so that paragraph's endParaRPr survives; it is stripped off at the very end
# with a targeted Delete() rather than clearing .Text (which would drop the
# paragraph mark formatting).
$tr.InsertBefore("X")
$tr.InsertBefore("Finding which data to focus on`r")
$tr.InsertBefore("Potentially put this section later, when talking about doing things differently in the future`r")
$tr.InsertBefore("Issues that came up`r")
$tr.InsertBefore("Screenshot or example of the data we found through external sources`r")

# Paragraph 1 - "Screenshot or example of the data we found through external sources"
$para1 = $tr.Paragraphs(1, 1)
$para1.Font.Language = 1033
$para1.IndentLevel = 1
$para1.ParagraphFormat.Bullet.Visible = $true
$para1.ParagraphFormat.Bullet.Font.Name = "Arial"
$para1.ParagraphFormat.Bullet.Character = 8226

# Paragraph 2 - "Issues that came up"
$para2 = $tr.Paragraphs(2, 1)
$para2.Font.Language = 1033
$para2.IndentLevel = 1
$para2.ParagraphFormat.Bullet.Visible = $true
$para2.ParagraphFormat.Bullet.Font.Name = "Arial"
$para2.ParagraphFormat.Bullet.Character = 8226

# Paragraph 3 - "Potentially put this section later, ..." (sub-level)
$para3 = $tr.Paragraphs(3, 1)
$para3.Font.Language = 1033
$para3.IndentLevel = 2
$para3.ParagraphFormat.Bullet.Visible = $true
$para3.ParagraphFormat.Bullet.Font.Name = "Arial"
$para3.ParagraphFormat.Bullet.Character = 8226

# Paragraph 4 - "Finding which data to " + "focus on" (two runs)
$para4 = $tr.Paragraphs(4, 1)
$para4.IndentLevel = 1
$para4.ParagraphFormat.Bullet.Visible = $true
$para4.ParagraphFormat.Bullet.Font.Name = "Arial"
$para4.ParagraphFormat.Bullet.Character = 8226
$para4.Text = "Finding which data to "
$para4.Font.Language = 1033
$para4.InsertAfter("focus on")

# Paragraph 5 - trailing empty paragraph (sub-level, no bullet)
$para5 = $tr.Paragraphs(5, 1)
$para5.IndentLevel = 2
$para5.ParagraphFormat.Bullet.Visible = $false
$para5.Delete()
